$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.877.33"
$ws.Range("E2").Value = "  +2.67%  "

$ws.Range("D3").Value = "1.668.79"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("D5").Value = "'214.78"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "'23.59"
$ws.Range("E8").Value = "  +3.59%  "

$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("E11").Value = "  -1.27%  "

$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").Value = "1.671.33"
$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("D14").Value = "'4.16"

$ws.Range("E15").Value = "  -1.32%  "

$ws.Range("D16").Value = "'66.08"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").Value = "'251.23"
$ws.Range("E17").Value = "  +6.84%  "

$ws.Range("D18").Value = "27.851.42"
$ws.Range("E18").Value = "  +2.69%  "

$ws.Range("E19").Value = "  -1.19%  "

$ws.Range("D20").Value = "'7.56"
$ws.Range("E20").Value = "  -3.87%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("E22").Value = "  -1.34%  "

$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("E24").Value = "  -1.58%  "

$ws.Range("D25").Value = "'146.90"
$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("E26").Value = "  -3.03%  "

$ws.Range("D27").Value = "'16.33"
$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("E32").Value = "  -0.45%  "

$ws.Range("E33").Value = "  -2.60%  "

$ws.Range("D34").Value = "1.430.71"
$ws.Range("E34").Value = "  -7.14%  "

$ws.Range("E35").Value = "  -5.72%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").Value = "'0.931"
$ws.Range("E37").Value = "  -1.38%  "

$ws.Range("D38").Value = "'0.584"
$ws.Range("E38").Value = "  -3.91%  "

$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "'1.04"
$ws.Range("E40").Value = "  -2.46%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'69.77"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.40"
$ws.Range("E43").Value = "  -6.56%  "

$ws.Range("D44").Value = "1.811.77"
$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "'2.21"
$ws.Range("E45").Value = "  -1.71%  "

$ws.Range("D46").Value = "'0.790"
$ws.Range("E46").Value = "  +1.23%  "

$ws.Range("E47").Value = "  +4.89%  "

$ws.Range("D48").Value = "'89.04"
$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -3.51%  "

$ws.Range("E50").Value = "  -1.95%  "

$ws.Range("D51").Value = "'7.85"
$ws.Range("E51").Value = "  -4.99%  "
